$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data refresh: values rescaled to 900 minutes and rows re-sorted by
# descending |Diff. Top 5 avec Bottom 15 en %| (column D).
$data = @{
  2 = @{
    A = "highdecel_count_full_tip"
    B = 40.32643280579968
    C = 37.81427744046337
    D = 6.64
    E = 1.566944369617723
    F = 3.277827349723662
    G = 37.69567533837247
    H = 33.918905803945
    I = 41.65753511656836
    J = 44.94408914964191
  }
  3 = @{
    A = "total_distance_full_otip"
    B = 35557.73874765011
    C = 37875.22498564087
    D = -6.12
    E = 1232.388864977056
    F = 2824.068090997378
    G = 34037.08639854103
    H = 31960.22600064877
    I = 37278.73176481328
    J = 41842.81071672353
  }
  4 = @{
    A = "sprint_distance_full_tip"
    B = 781.8154777733826
    C = 742.2055356201979
    D = 5.34
    E = 46.33414738302128
    F = 71.70087241178432
    G = 713.0034047802379
    H = 621.2318161859457
    I = 840.5826305662382
    J = 856.2483521944293
  }
  5 = @{
    A = "highaccel_count_full_tip"
    B = 21.5808006642409
    C = 20.58435855497999
    D = 4.84
    E = 2.075893813909885
    F = 2.484906787937065
    G = 18.66445279530906
    H = 17.21672867397698
    I = 24.12555135755219
    J = 26.38165588884013
  }
  6 = @{
    A = "meddecel_count_full_otip"
    B = 295.7264310572242
    C = 310.6382529620406
    D = -4.8
    E = 21.48979728892436
    F = 26.68113240161935
    G = 271.3609392051098
    H = 250.370317805314
    I = 328.8671087604378
    J = 352.7150321250502
  }
  7 = @{
    A = "meddecel_count_full_tip"
    B = 260.9359912541216
    C = 249.8672855776814
    D = 4.43
    E = 16.95666447451088
    F = 18.5851248816636
    G = 233.8347715843163
    H = 223.6140481905484
    I = 278.0250160352341
    J = 289.7676027772649
  }
  8 = @{
    A = "total_distance_full_tip"
    B = 35483.45315552852
    C = 34213.03794152191
    D = 3.71
    E = 1779.777983213065
    F = 2197.813938713118
    G = 32647.87026909309
    H = 30965.54675942479
    I = 36984.78942497371
    J = 38666.44617415352
  }
  9 = @{
    A = "sprint_count_full_tip"
    B = 39.34931314487949
    C = 37.94439031012171
    D = 3.7
    E = 2.27103105593123
    F = 3.672266351188446
    G = 36.12928511257346
    H = 32.77406582441426
    I = 42.10090083547875
    J = 44.43655945661627
  }
  10 = @{
    A = "running_distance_full_otip"
    B = 6853.733652790443
    C = 7090.395590995666
    D = -3.34
    E = 261.7312196234597
    F = 824.1695915250681
    G = 6461.814865793849
    H = 5857.257841365194
    I = 7116.883714345494
    J = 9202.918048062398
  }
  11 = @{
    A = "medaccel_count_full_tip"
    B = 367.4935104648872
    C = 355.7742733909068
    D = 3.29
    E = 26.73875633393764
    F = 23.04717293601096
    G = 325.3385808165229
    H = 323.6035574796436
    I = 393.013095416546
    J = 403.1707345705044
  }
  12 = @{
    A = "running_distance_full_tip"
    B = 5652.00052074328
    C = 5477.507365740687
    D = 3.19
    E = 398.953218813526
    F = 388.379698055164
    G = 5243.529566463166
    H = 4769.842932808218
    I = 6237.721800697452
    J = 6193.01590876589
  }
  13 = @{
    A = "medaccel_count_full_otip"
    B = 413.3315748917781
    C = 425.1069623861311
    D = -2.77
    E = 19.46897222949351
    F = 30.51045719714014
    G = 396.2543559713862
    H = 360.4221816733501
    I = 444.7572931794446
    J = 467.8653562805059
  }
  14 = @{
    A = "sprint_distance_full_otip"
    B = 707.7137791400927
    C = 727.1140904712177
    D = -2.67
    E = 37.16671986120333
    F = 125.5787783588835
    G = 654.4873509197188
    H = 557.7131031458366
    I = 757.7201179376283
    J = 1100.605823962163
  }
  15 = @{
    A = "sprint_count_full_otip"
    B = 37.16239915648865
    C = 38.1755827592341
    D = -2.65
    E = 1.85777303925926
    F = 7.204822822810137
    G = 34.32881547899326
    H = 29.90705738630099
    I = 39.44895443994156
    J = 60.41563664314325
  }
  16 = @{
    A = "hi_distance_full_tip"
    B = 2983.627093050226
    C = 2906.84415604861
    D = 2.64
    E = 192.6248345487242
    F = 228.6704536234355
    G = 2741.886148356625
    H = 2574.649203297054
    I = 3187.861066187726
    J = 3292.476696441771
  }
  17 = @{
    A = "hi_count_full_tip"
    B = 248.7348905139809
    C = 242.7086526873697
    D = 2.48
    E = 17.97719918149894
    F = 18.18235700935413
    G = 225.1166050882141
    H = 216.4225409928567
    I = 267.7479426716968
    J = 276.4090780947945
  }
  18 = @{
    A = "hsr_count_full_tip"
    B = 209.3855773691015
    C = 204.764262377248
    D = 2.26
    E = 15.76290027667297
    F = 14.67104467656554
    G = 188.9873199756406
    H = 183.1214673779982
    I = 225.6470418362181
    J = 231.9725186381782
  }
  19 = @{
    A = "highaccel_count_full_otip"
    B = 20.9885927830018
    C = 20.55151676233293
    D = 2.13
    E = 1.644643387792927
    F = 1.072954953754047
    G = 19.01366400801037
    H = 19.16336104346119
    I = 23.50821093050978
    J = 22.64911874439506
  }
  20 = @{
    A = "hsr_distance_full_tip"
    B = 2201.811615276843
    C = 2164.638620428413
    D = 1.72
    E = 162.1337972889364
    F = 165.4497535123431
    G = 2028.882743576388
    H = 1895.187098614456
    I = 2404.758738272372
    J = 2453.299623878735
  }
  21 = @{
    A = "hi_count_full_otip"
    B = 285.5495749623238
    C = 290.2838287275494
    D = -1.63
    E = 14.87682812721861
    F = 41.51589737945454
    G = 262.0781747488301
    H = 236.0281904784996
    I = 301.9106687576135
    J = 413.3904060305224
  }
  22 = @{
    A = "hsr_count_full_otip"
    B = 248.3871758058353
    C = 252.1082459683153
    D = -1.48
    E = 13.07014610713024
    F = 34.43216839108684
    G = 227.7493592698368
    H = 206.092053379901
    I = 262.4617143176719
    J = 352.9747693873791
  }
  23 = @{
    A = "highdecel_count_full_otip"
    B = 51.72276474536343
    C = 52.32048826977314
    D = -1.14
    E = 4.036793686977469
    F = 5.826276827397506
    G = 47.46093844235369
    H = 44.40391729390301
    I = 56.52843416597797
    J = 67.31286071329997
  }
  24 = @{
    A = "hi_distance_full_otip"
    B = 3413.389484521112
    C = 3435.329989578204
    D = -0.64
    E = 173.7964389024005
    F = 556.3293553007516
    G = 3166.75182288359
    H = 2803.614667939856
    I = 3656.218461286545
    J = 5093.182468857197
  }
  25 = @{
    A = "total_metersperminute_full_otip"
    B = 1921.369056708694
    C = 1933.397142085327
    D = -0.62
    E = 117.5788256070493
    F = 75.84925843295528
    G = 1741.475776568035
    H = 1812.256080771003
    I = 2055.524019744857
    J = 2147.247733113258
  }
  26 = @{
    A = "total_metersperminute_full_tip"
    B = 1807.588149029572
    C = 1812.821493985028
    D = -0.29
    E = 110.6389375267606
    F = 41.61203321987711
    G = 1690.440770990897
    H = 1752.370680344522
    I = 1971.081921566278
    J = 1886.173303597148
  }
  27 = @{
    A = "hsr_distance_full_otip"
    B = 2705.675705381019
    C = 2708.215899106986
    D = -0.09
    E = 137.4437812126859
    F = 432.8025992886028
    G = 2512.264471963871
    H = 2240.149115371444
    I = 2898.498343348916
    J = 3992.576644895034
  }
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  $ws.Cells.Item([int]$r, 1).Value = $row.A
  $ws.Cells.Item([int]$r, 2).Value = $row.B
  $ws.Cells.Item([int]$r, 3).Value = $row.C
  $ws.Cells.Item([int]$r, 4).Value = $row.D
  $ws.Cells.Item([int]$r, 5).Value = $row.E
  $ws.Cells.Item([int]$r, 6).Value = $row.F
  $ws.Cells.Item([int]$r, 7).Value = $row.G
  $ws.Cells.Item([int]$r, 8).Value = $row.H
  $ws.Cells.Item([int]$r, 9).Value = $row.I
  $ws.Cells.Item([int]$r, 10).Value = $row.J
}
